$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Halve the G9:G38 values (new passenger car demand figures)
for ($r = 9; $r -le 38; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $old = $cell.Value()
    $cell.Value = $old / 2
}

# Update selection to reflect plotting range G46:G78
$ws.Range("G46:G78").Select()
